$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, pushing the existing row 13 (and below) down to row 14.
$ws.Rows.Item(13).Insert()

# The old row 13 data now lives at row 14 (Excel copies formatting down but not values,
# so we explicitly restore the original row-13 values into row 14).
$ws.Cells.Item(14, 1).Value = 1
$ws.Cells.Item(14, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(14, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(14, 4).Value = 44281
$ws.Cells.Item(14, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(14, 5).Value = 15
$ws.Cells.Item(14, 6).Value = 100112017
$ws.Cells.Item(14, 7).Value = "Ramas de apio"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 100
$ws.Cells.Item(14, 11).Value = 5000
$ws.Cells.Item(14, 12).Value = 6000
$ws.Cells.Item(14, 13).Value = 5500
$ws.Cells.Item(14, 14).Value = "$/atado 7 kilos"
$ws.Cells.Item(14, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(14, 16).Value = 5500
$ws.Cells.Item(14, 17).Value = 1
$ws.Cells.Item(14, 18).Value = "Hortaliza"

# New row 13 gets the fresh weekly entry.
$ws.Cells.Item(13, 1).Value = 1
$ws.Cells.Item(13, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(13, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(13, 4).Value = 44636
$ws.Cells.Item(13, 5).Value = 15
$ws.Cells.Item(13, 6).Value = 100112017
$ws.Cells.Item(13, 7).Value = "Ramas de apio"
$ws.Cells.Item(13, 8).Value = "Americana (o)"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 60
$ws.Cells.Item(13, 11).Value = 8000
$ws.Cells.Item(13, 12).Value = 9000
$ws.Cells.Item(13, 13).Value = 8500
$ws.Cells.Item(13, 14).Value = "$/atado 7 kilos"
$ws.Cells.Item(13, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(13, 16).Value = 8500
$ws.Cells.Item(13, 17).Value = 1
$ws.Cells.Item(13, 18).Value = "Hortaliza"
